$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Battery_MV")

$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 24000

$ws.Range("G2").Formula = "=1000"
$ws.Range("G2").HorizontalAlignment = -4108

$ws.Range("H2").Value = 24000
$ws.Range("H2").HorizontalAlignment = -4108

$ws.Range("I2").Value = 1000

$ws.Range("L2").ClearFormats()
$ws.Range("L2").Formula = "=-1/0.2"

$ws.Activate()
$ws.Range("D8").Select()
